# Update "想去人数" (F column) counts for both the "展览" and "全部类型" sheets
# to reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new F value for the "展览" sheet
$exhibitionUpdates = @{
    4  = 8043
    5  = 99
    8  = 32
    10 = 479
    13 = 459
    14 = 71
    15 = 80
    17 = 5943
    18 = 191
    19 = 282
    20 = 1990
    21 = 43
    22 = 64
    23 = 240
    24 = 414
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Row -> new F value for the "全部类型" sheet
$allTypesUpdates = @{
    4  = 8043
    5  = 99
    8  = 32
    10 = 479
    13 = 459
    14 = 71
    15 = 80
    18 = 5943
    20 = 191
    21 = 282
    22 = 1990
    23 = 43
    24 = 64
    25 = 240
    26 = 414
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}

$wb.Save()
